$d = $word.ActiveDocument

# --- Edit 1: intro sentence above the list of top-performing schools ---
# "school have ... They are in order of highest to lowest."
#   -> "schools have ... Displayed below in descending order."
$old1 = "The below school have the highest combine passing rate of 70% or above.  They are in order of highest to lowest. "
$new1 = "The below schools have the highest combine passing rate of 70% or above.  Displayed below in descending order. "
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- Edit 2: the "_GoBack" bookmark currently sits at the end of the
# "Written 2" paragraph; it needs to move to the end of the following
# paragraph (just before the final period), once that paragraph's text
# is rewritten below. Remove it from its old spot first. ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Edit 3: rewrite the "low performing schools" paragraph ---
$old2 = "The low performing schools have a greater dollars per student spend than the top five performing school."
$new2 = "The lower performing schools have spent more  per student spend than the top five performing schools."
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# --- Edit 4: re-insert the "_GoBack" bookmark right before the trailing
# period of that rewritten paragraph (collapsed/zero-length, as before) ---
$locate = $d.Range(0, $d.Content.End)
$locate.Find.Execute("performing schools.") | Out-Null
$periodPos = $locate.End - 1
$bmRange = $d.Range($periodPos, $periodPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
